# Update NATMI TPM-derived LR-pair metrics with recalculated TPM values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"25.35940266666667"
$ws.Range("H2").Value = [double]"76.078208"
$ws.Range("I2").Value = [double]"0.005186643687654987"
$ws.Range("J2").Value = [double]"0.005186643687654986"
$ws.Range("M2").Value = [double]"3.319447"
$ws.Range("N2").Value = [double]"9.958341000000001"
$ws.Range("O2").Value = [double]"0.03276237985858125"
$ws.Range("P2").Value = [double]"0.03276237985858125"
$ws.Range("Q2").Value = [double]"84.17919310365868"
$ws.Range("R2").Value = [double]"757.6127379329281"
$ws.Range("S2").Value = [double]"0.0001699267906860653"
$ws.Range("T2").Value = [double]"0.0001699267906860653"
$ws.Range("G3").Value = [double]"25.35940266666667"
$ws.Range("H3").Value = [double]"76.078208"
$ws.Range("I3").Value = [double]"0.005186643687654987"
$ws.Range("J3").Value = [double]"0.005186643687654986"
$ws.Range("O3").Value = [double]"0.7010975337260504"
$ws.Range("P3").Value = [double]"0.7010975337260504"
$ws.Range("Q3").Value = [double]"1801.390037316409"
$ws.Range("R3").Value = [double]"16212.51033584768"
$ws.Range("S3").Value = [double]"0.003636343097730699"
$ws.Range("T3").Value = [double]"0.003636343097730698"
$ws.Range("G4").Value = [double]"25.35940266666667"
$ws.Range("H4").Value = [double]"76.078208"
$ws.Range("I4").Value = [double]"0.005186643687654987"
$ws.Range("J4").Value = [double]"0.005186643687654986"
$ws.Range("M4").Value = [double]"26.57769466666667"
$ws.Range("N4").Value = [double]"79.73308400000001"
$ws.Range("O4").Value = [double]"0.262317346363633"
$ws.Range("P4").Value = [double]"0.262317346363633"
$ws.Range("Q4").Value = [double]"673.9944610037193"
$ws.Range("R4").Value = [double]"6065.950149033472"
$ws.Range("S4").Value = [double]"0.001360546608679344"
$ws.Range("T4").Value = [double]"0.001360546608679343"
$ws.Range("G5").Value = [double]"25.35940266666667"
$ws.Range("H5").Value = [double]"76.078208"
$ws.Range("I5").Value = [double]"0.005186643687654987"
$ws.Range("J5").Value = [double]"0.005186643687654986"
$ws.Range("M5").Value = [double]"0.3873156666666667"
$ws.Range("N5").Value = [double]"1.161947"
$ws.Range("O5").Value = [double]"0.003822740051735415"
$ws.Range("P5").Value = [double]"0.003822740051735415"
$ws.Range("Q5").Value = [double]"9.822093950108444"
$ws.Range("R5").Value = [double]"88.39884555097601"
$ws.Range("S5").Value = [double]"1.982719055887939E-05"
$ws.Range("T5").Value = [double]"1.982719055887939E-05"
$ws.Range("I6").Value = [double]"0.9837462940761621"
$ws.Range("J6").Value = [double]"0.983746294076162"
$ws.Range("M6").Value = [double]"3.319447"
$ws.Range("N6").Value = [double]"9.958341000000001"
$ws.Range("O6").Value = [double]"0.03276237985858125"
$ws.Range("P6").Value = [double]"0.03276237985858125"
$ws.Range("Q6").Value = [double]"15966.19591416097"
$ws.Range("R6").Value = [double]"143695.7632274487"
$ws.Range("S6").Value = [double]"0.0322298697709948"
$ws.Range("T6").Value = [double]"0.0322298697709948"
$ws.Range("I7").Value = [double]"0.9837462940761621"
$ws.Range("J7").Value = [double]"0.983746294076162"
$ws.Range("O7").Value = [double]"0.7010975337260504"
$ws.Range("P7").Value = [double]"0.7010975337260504"
$ws.Range("S7").Value = [double]"0.6897021005889392"
$ws.Range("T7").Value = [double]"0.6897021005889391"
$ws.Range("I8").Value = [double]"0.9837462940761621"
$ws.Range("J8").Value = [double]"0.983746294076162"
$ws.Range("M8").Value = [double]"26.57769466666667"
$ws.Range("N8").Value = [double]"79.73308400000001"
$ws.Range("O8").Value = [double]"0.262317346363633"
$ws.Range("P8").Value = [double]"0.262317346363633"
$ws.Range("Q8").Value = [double]"127835.9558067206"
$ws.Range("R8").Value = [double]"1150523.602260485"
$ws.Range("S8").Value = [double]"0.258053717357117"
$ws.Range("T8").Value = [double]"0.2580537173571169"
$ws.Range("I9").Value = [double]"0.9837462940761621"
$ws.Range("J9").Value = [double]"0.983746294076162"
$ws.Range("M9").Value = [double]"0.3873156666666667"
$ws.Range("N9").Value = [double]"1.161947"
$ws.Range("O9").Value = [double]"0.003822740051735415"
$ws.Range("P9").Value = [double]"0.003822740051735415"
$ws.Range("Q9").Value = [double]"1862.948200294768"
$ws.Range("R9").Value = [double]"16766.53380265291"
$ws.Range("S9").Value = [double]"0.00376060635911123"
$ws.Range("T9").Value = [double]"0.00376060635911123"
$ws.Range("G10").Value = [double]"51.27300266666666"
$ws.Range("H10").Value = [double]"153.819008"
$ws.Range("I10").Value = [double]"0.01048663484403512"
$ws.Range("J10").Value = [double]"0.01048663484403512"
$ws.Range("M10").Value = [double]"3.319447"
$ws.Range("N10").Value = [double]"9.958341000000001"
$ws.Range("O10").Value = [double]"0.03276237985858125"
$ws.Range("P10").Value = [double]"0.03276237985858125"
$ws.Range("Q10").Value = [double]"170.1980148828587"
$ws.Range("R10").Value = [double]"1531.782133945728"
$ws.Range("S10").Value = [double]"0.0003435671141985127"
$ws.Range("T10").Value = [double]"0.0003435671141985127"
$ws.Range("G11").Value = [double]"51.27300266666666"
$ws.Range("H11").Value = [double]"153.819008"
$ws.Range("I11").Value = [double]"0.01048663484403512"
$ws.Range("J11").Value = [double]"0.01048663484403512"
$ws.Range("O11").Value = [double]"0.7010975337260504"
$ws.Range("P11").Value = [double]"0.7010975337260504"
$ws.Range("Q11").Value = [double]"3642.147151535076"
$ws.Range("R11").Value = [double]"32779.32436381569"
$ws.Range("S11").Value = [double]"0.007352153826238691"
$ws.Range("T11").Value = [double]"0.007352153826238691"
$ws.Range("G12").Value = [double]"51.27300266666666"
$ws.Range("H12").Value = [double]"153.819008"
$ws.Range("I12").Value = [double]"0.01048663484403512"
$ws.Range("J12").Value = [double]"0.01048663484403512"
$ws.Range("M12").Value = [double]"26.57769466666667"
$ws.Range("N12").Value = [double]"79.73308400000001"
$ws.Range("O12").Value = [double]"0.262317346363633"
$ws.Range("P12").Value = [double]"0.262317346363633"
$ws.Range("Q12").Value = [double]"1362.718209517853"
$ws.Range("R12").Value = [double]"12264.46388566067"
$ws.Range("S12").Value = [double]"0.002750826224571704"
$ws.Range("T12").Value = [double]"0.002750826224571704"
$ws.Range("G13").Value = [double]"51.27300266666666"
$ws.Range("H13").Value = [double]"153.819008"
$ws.Range("I13").Value = [double]"0.01048663484403512"
$ws.Range("J13").Value = [double]"0.01048663484403512"
$ws.Range("M13").Value = [double]"0.3873156666666667"
$ws.Range("N13").Value = [double]"1.161947"
$ws.Range("O13").Value = [double]"0.003822740051735415"
$ws.Range("P13").Value = [double]"0.003822740051735415"
$ws.Range("Q13").Value = [double]"19.85883720984178"
$ws.Range("R13").Value = [double]"178.729534888576"
$ws.Range("S13").Value = [double]"4.008767902621724E-05"
$ws.Range("T13").Value = [double]"4.008767902621724E-05"
$ws.Range("G14").Value = [double]"2.837922333333333"
$ws.Range("H14").Value = [double]"8.513767"
$ws.Range("I14").Value = [double]"0.0005804273921477663"
$ws.Range("J14").Value = [double]"0.0005804273921477662"
$ws.Range("M14").Value = [double]"3.319447"
$ws.Range("N14").Value = [double]"9.958341000000001"
$ws.Range("O14").Value = [double]"0.03276237985858125"
$ws.Range("P14").Value = [double]"0.03276237985858125"
$ws.Range("Q14").Value = [double]"9.420332775616334"
$ws.Range("R14").Value = [double]"84.782994980547"
$ws.Range("S14").Value = [double]"1.901618270187082E-05"
$ws.Range("T14").Value = [double]"1.901618270187082E-05"
$ws.Range("G15").Value = [double]"2.837922333333333"
$ws.Range("H15").Value = [double]"8.513767"
$ws.Range("I15").Value = [double]"0.0005804273921477663"
$ws.Range("J15").Value = [double]"0.0005804273921477662"
$ws.Range("O15").Value = [double]"0.7010975337260504"
$ws.Range("P15").Value = [double]"0.7010975337260504"
$ws.Range("Q15").Value = [double]"201.5901196546745"
$ws.Range("R15").Value = [double]"1814.31107689207"
$ws.Range("S15").Value = [double]"0.0004069362131418421"
$ws.Range("T15").Value = [double]"0.000406936213141842"
$ws.Range("G16").Value = [double]"2.837922333333333"
$ws.Range("H16").Value = [double]"8.513767"
$ws.Range("I16").Value = [double]"0.0005804273921477663"
$ws.Range("J16").Value = [double]"0.0005804273921477662"
$ws.Range("M16").Value = [double]"26.57769466666667"
$ws.Range("N16").Value = [double]"79.73308400000001"
$ws.Range("O16").Value = [double]"0.262317346363633"
$ws.Range("P16").Value = [double]"0.262317346363633"
$ws.Range("Q16").Value = [double]"75.42543326304757"
$ws.Range("R16").Value = [double]"678.828899367428"
$ws.Range("S16").Value = [double]"0.0001522561732649658"
$ws.Range("T16").Value = [double]"0.0001522561732649658"
$ws.Range("G17").Value = [double]"2.837922333333333"
$ws.Range("H17").Value = [double]"8.513767"
$ws.Range("I17").Value = [double]"0.0005804273921477663"
$ws.Range("J17").Value = [double]"0.0005804273921477662"
$ws.Range("M17").Value = [double]"0.3873156666666667"
$ws.Range("N17").Value = [double]"1.161947"
$ws.Range("O17").Value = [double]"0.003822740051735415"
$ws.Range("P17").Value = [double]"0.003822740051735415"
$ws.Range("Q17").Value = [double]"1.099171780483222"
$ws.Range("R17").Value = [double]"9.892546024349"
$ws.Range("S17").Value = [double]"2.218823039087604E-06"
$ws.Range("T17").Value = [double]"2.218823039087604E-06"
